# Updates the cached "datetimeFigureOut" date field text (8/8/23 -> 10/5/23)
# on every slide layout, the slide master and the notes master, and adds a
# yellow highlight to the "X regions" / "Y regions" runs on slide 1.

$p = $ppt.ActivePresentation

$oldDate = "8/8/23"
$newDate = "10/5/23"

function Update-DatePlaceholder {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# 1. Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# 2. Every slide layout hanging off the master
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# 3. Notes master
$notesMaster = $p.NotesMaster
Update-DatePlaceholder $notesMaster.Shapes

# 4. Highlight "X regions" / "Y regions" runs in yellow on slide 1
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        $fullText = $tr.Text
        foreach ($label in @("X regions", "Y regions")) {
            $pos = $fullText.IndexOf($label)
            if ($pos -ge 0) {
                $sub = $tr.Characters($pos + 1, $label.Length)
                $sub.Font.Highlight = 65535
            }
        }
    }
}
